$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (exhibitions) - F-column "want to go" counters bump up, and
# G5 flips from a numeric price to the text "不可售" (not sellable).
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

$wsExpo.Range("F2").Value  = 33
$wsExpo.Range("F3").Value  = 21135
$wsExpo.Range("F4").Value  = 816
$wsExpo.Range("G5").Value  = "不可售"
$wsExpo.Range("F6").Value  = 1126
$wsExpo.Range("F7").Value  = 25
$wsExpo.Range("F8").Value  = 7884
$wsExpo.Range("F10").Value = 40
$wsExpo.Range("F12").Value = 306
$wsExpo.Range("F15").Value = 163
$wsExpo.Range("F18").Value = 224
$wsExpo.Range("F19").Value = 1357
$wsExpo.Range("F20").Value = 516
$wsExpo.Range("F27").Value = 1176
$wsExpo.Range("F28").Value = 50
$wsExpo.Range("F29").Value = 38
$wsExpo.Range("F32").Value = 599
$wsExpo.Range("F34").Value = 133
$wsExpo.Range("F35").Value = 5022
$wsExpo.Range("F40").Value = 13040
$wsExpo.Range("F44").Value = 69
$wsExpo.Range("F46").Value = 422
$wsExpo.Range("F47").Value = 4056
$wsExpo.Range("F48").Value = 330
$wsExpo.Range("F49").Value = 101

# ---------------------------------------------------------------------------
# Sheet "全部类型" (all types) - same kind of F-column bumps, plus row 4
# is replaced: the old "coke老师撸猫内场票" listing is swapped out for the
# new "OCG国潮动漫游戏嘉年华·你的欲梦内场" listing (already present on the
# "展览" sheet).
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Range("F2").Value = 33
$wsAll.Range("F3").Value = 21135

$wsAll.Range("C4").Value = "【大会员提前抢】苏州·OCG国潮动漫游戏嘉年华·你的欲梦内场·全网内场首签"
$wsAll.Range("D4").Value = "苏州大道东688号 苏州国际博览中心"
$wsAll.Range("E4").Value = "2024.07.20 09:00-07.20 17:00"
$wsAll.Range("F4").Value = 816
$wsAll.Range("G4").Value = "已售罄"
$wsAll.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=86884"
$wsAll.Range("I4").Value = "//i1.hdslb.com/bfs/openplatform/202406/Ir2IFJnI1717538486594.jpeg"

$wsAll.Range("F5").Value  = 1126
$wsAll.Range("F6").Value  = 25
$wsAll.Range("F7").Value  = 7884
$wsAll.Range("F9").Value  = 40
$wsAll.Range("F11").Value = 306
$wsAll.Range("F14").Value = 163
$wsAll.Range("F16").Value = 224
$wsAll.Range("F17").Value = 1357
$wsAll.Range("F18").Value = 516
$wsAll.Range("F25").Value = 1176
$wsAll.Range("F26").Value = 50
$wsAll.Range("F27").Value = 38
$wsAll.Range("F30").Value = 599
$wsAll.Range("F33").Value = 133
$wsAll.Range("F35").Value = 5022
$wsAll.Range("F40").Value = 13040
$wsAll.Range("F44").Value = 69
$wsAll.Range("F46").Value = 422
$wsAll.Range("F47").Value = 4056
$wsAll.Range("F48").Value = 330
$wsAll.Range("F49").Value = 101

$wb.Save()
